# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the Status column moves
# from "Handed back: in sync with en-US" to "Ready for handoff", the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# are refreshed, and the (now much shorter) Status column is re-sized to
# fit its new, shorter content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed timestamps ---
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff
# Datetime" shared the same old value (2016-08-15 16:55:14); both move to
# the new handoff-generation timestamp.
$wsOverview.Range("G2").Value = "2016-08-15 16:55:59"
$wsDeDe.Range("H2").Value     = "2016-08-15 16:55:59"

# zh-cn!H2 "Latest Handoff Datetime" gets its own refreshed timestamp.
$wsZhCn.Range("H2").Value = "2016-08-15 16:55:54"

# --- Column width: the Status column shrinks to fit the shorter text ---
# The report's own column auto-sizer targets ~17.22 "characters"; the
# Excel COM ColumnWidth setter only resolves to the nearest 1/6th of a
# character (pixel-snapped), so 16.3333... is the input that lands closest
# on the intended width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3333333333333
